$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 401, shifting the existing rows
# (401-417) down to (403-419).
$ws.Rows.Item(401).Insert()
$ws.Rows.Item(401).Insert()

# --- New row 401: June Pearl / Primera ---
$ws.Range("A401").Value = 4
$ws.Range("B401").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C401").Value = "Los Lagos"
$ws.Range("D401").Value = 44610
$ws.Range("E401").Value = 10
$ws.Range("F401").Value = "Fruta"
$ws.Range("G401").Value = 100103
$ws.Range("H401").Value = "Frutos de hueso (carozo)"
$ws.Range("I401").Value = 100103006
$ws.Range("J401").Value = "Nectarín"
$ws.Range("K401").Value = "June Pearl"
$ws.Range("L401").Value = "Primera"
$ws.Range("M401").Value = 500
$ws.Range("N401").Value = 16000
$ws.Range("O401").Value = 16000
$ws.Range("P401").Value = 16000
$ws.Range("Q401").Value = "$/caja 15 kilos empedrada"
$ws.Range("R401").Value = "Región de O'Higgins"
$ws.Range("S401").Value = 1067
$ws.Range("T401").Value = 15

# --- New row 402: Venus / Primera ---
$ws.Range("A402").Value = 4
$ws.Range("B402").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C402").Value = "Los Lagos"
$ws.Range("D402").Value = 44610
$ws.Range("E402").Value = 10
$ws.Range("F402").Value = "Fruta"
$ws.Range("G402").Value = 100103
$ws.Range("H402").Value = "Frutos de hueso (carozo)"
$ws.Range("I402").Value = 100103006
$ws.Range("J402").Value = "Nectarín"
$ws.Range("K402").Value = "Venus"
$ws.Range("L402").Value = "Primera"
$ws.Range("M402").Value = 500
$ws.Range("N402").Value = 16000
$ws.Range("O402").Value = 16000
$ws.Range("P402").Value = 16000
$ws.Range("Q402").Value = "$/caja 15 kilos empedrada"
$ws.Range("R402").Value = "Región de O'Higgins"
$ws.Range("S402").Value = 1067
$ws.Range("T402").Value = 15
